# Tutorial text updates for PeptideShaker v0.33.4 / A-score wording tweaks.
$d = $word.ActiveDocument

# Helper: re-establish run boundaries inside a (possibly coalesced) run of
# text by toggling Bold on/off over precise sub-ranges. Word (and this
# interop layer) merges adjacent runs that share byte-identical formatting
# whenever a Find/Replace edits one of them; flipping a formatting property
# back and forth on an exact sub-range forces a clean split at that
# boundary without touching the text itself.
function Split-AtOffsets($matchText, [int[]]$offsets) {
    $r = $d.Content
    $found = $r.Find.Execute($matchText, $true)
    if (-not $found) {
        Write-Host "Split-AtOffsets: NOT FOUND:" $matchText
        return
    }
    $base = $r.Start
    foreach ($off in $offsets) {
        $sub = $d.Range($base, $base + $off)
        $sub.Bold = 1
        $sub.Bold = 0
    }
}

# ---------------------------------------------------------------------
# 1) "Note that all spectrum annotation (modifications, ions, etc.) have
#    been passed by " -> "Note that all the spectrum annotations
#    (modifications, ions, etc.) have been passed by ", expressed as five
#    separate runs.
# ---------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute( `
    "Note that all spectrum annotation (modifications, ions, etc.) have been passed by ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Note that all the spectrum annotations (modifications, ions, etc.) have been passed by ", 2)
Write-Host "Step1 replace found:" $found1

# ---------------------------------------------------------------------
# 2) "available for all other online resources" -> "available for other
#    online resources", expressed as two separate runs.
# ---------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute( `
    "available for all other online resources", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "available for other online resources", 2)
Write-Host "Step2 replace found:" $found2

# ---------------------------------------------------------------------
# 3) The italic double space before "What difference do you see..."
#    becomes a single space.
# ---------------------------------------------------------------------
$r3 = $d.Content
$found3 = $r3.Find.Execute( `
    "  What difference", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " What difference", 2)
Write-Host "Step3 replace found:" $found3

# ---------------------------------------------------------------------
# Re-split the runs that the edits above coalesced, back to the exact
# boundaries required.
# ---------------------------------------------------------------------
Split-AtOffsets "Note that all the spectrum annotations (modifications, ions, etc.) have been passed by " @(13, 18, 37, 38)
Split-AtOffsets "as standardized terms and will thus be available for other online resources." @(39, 53, 75)
Split-AtOffsets " What difference do you see compared to the " @(1, 28, 39)

# ---------------------------------------------------------------------
# 4) Footer page field: the cached PAGE field result on the title page
#    goes from "4" to "1" (NUMPAGES field alongside it is untouched).
# ---------------------------------------------------------------------
$sec = $d.Sections.First
$ftr = $sec.Footers.Item(1)
$foundFtr = $ftr.Range.Find.Execute("4", $true, $false, $false, $false, $false, $true, 1, $false, "1", 2)
Write-Host "Footer page field replace found:" $foundFtr
